$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-11 (A:D) hold the existing "Subsystem / Target / Progress / Due date"
# table rows; row 11 is the last one ("Robot movement"). Copy its formatting
# down to the three new rows (12-14) we are about to add so the new rows pick
# up the same borders/number formats/fonts used throughout the table.
$ws.Range("A11:D11").Copy()
$ws.Range("A12:D14").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Row 12: All / Fully integrate / (in progress) / 18 Mar 2018
$ws.Range("A12").Value = "All"
$ws.Range("B12").Value = "Fully integrate"
$ws.Range("C12").Style = "Neutral"
$ws.Range("D12").Value = 43177

# Row 13: All / Begin whole system testing / (in progress) / 18 Mar 2018
$ws.Range("A13").Value = "All"
$ws.Range("B13").Value = "Begin whole system testing"
$ws.Range("C13").Style = "Neutral"
$ws.Range("D13").Value = 43177

# Row 14: All / Complete project / (in progress) / 20 Mar 2018
$ws.Range("A14").Value = "All"
$ws.Range("B14").Value = "Complete project"
$ws.Range("C14").Style = "Neutral"
$ws.Range("D14").Value = 43179

# Move the active selection to where it ended up after the edits.
$ws.Range("F5").Select()
